$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
    @(5, 5),
    @(9, 9),
    @(5, 6),
    @(7, 7),
    @(7, 7),
    @(11, 11),
    @(8, 8),
    @(9, 9),
    @(5, 5),
    @(7, 7),
    @(7, 7),
    @(6, 6),
    @(8, 8),
    @(9, 9),
    @(5, 6),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(6, 7),
    @(5, 5),
    @(7, 7),
    @(4, 4),
    @(7, 7),
    @(5, 5),
    @(6, 6),
    @(8, 8),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
